$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3695
$ws.Range("E2").Value = 365
$ws.Range("F2").Value = 348
$ws.Range("G2").Value = 194
$ws.Range("H2").Value = 145
$ws.Range("I2").Value = -1
$ws.Range("J2").Value = 146
$ws.Range("K2").Value = 4811
$ws.Range("L2").Value = 3159
$ws.Range("M2").Value = 1652
$ws.Range("N2").Value = 900
$ws.Range("O2").Value = 752
$ws.Range("P2").Value = 330
$ws.Range("Q2").Value = 546
$ws.Range("R2").Value = -268
$ws.Range("S2").Value = -335
$ws.Range("T2").Value = 250
$ws.Range("U2").Value = 296
$ws.Range("V2").Value = 1769
$ws.Range("W2").Value = 9.869999999999999
$ws.Range("X2").Value = 3.93
$ws.Range("Y2").Value = -0.07000000000000001
$ws.Range("Z2").Value = 2.94
$ws.Range("AA2").Value = 191.26
$ws.Range("AB2").Value = 185.58
$ws.Range("AC2").Value = -8
$ws.Range("AD2").Value = -3474.18
$ws.Range("AE2").Value = 11952
$ws.Range("AF2").Value = 2.4
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AJ2").Value = 7677460
$ws.Range("AI2").ClearContents()

# Row 3
$ws.Range("D3").Value = 3058
$ws.Range("E3").Value = 329
$ws.Range("F3").Value = 314
$ws.Range("G3").Value = 298
$ws.Range("H3").Value = 196
$ws.Range("I3").Value = 114
$ws.Range("J3").Value = 82
$ws.Range("K3").Value = 4560
$ws.Range("L3").Value = 3777
$ws.Range("M3").Value = 782
$ws.Range("N3").Value = 624
$ws.Range("O3").Value = 158
$ws.Range("P3").Value = 570
$ws.Range("Q3").Value = 235
$ws.Range("R3").Value = -153
$ws.Range("S3").Value = -261
$ws.Range("T3").Value = 162
$ws.Range("U3").Value = 72
$ws.Range("V3").Value = 2459
$ws.Range("W3").Value = 10.76
$ws.Range("X3").Value = 6.42
$ws.Range("Y3").Value = 14.99
$ws.Range("Z3").Value = 4.19
$ws.Range("AA3").Value = 482.87
$ws.Range("AB3").Value = -3.06
$ws.Range("AC3").Value = 1207
$ws.Range("AD3").Value = 44.19
$ws.Range("AE3").Value = 5231
$ws.Range("AF3").Value = 10.2
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AJ3").Value = 12088268
$ws.Range("AI3").ClearContents()

# Row 4
$ws.Range("D4").Value = 3352
$ws.Range("E4").Value = 289
$ws.Range("F4").Value = 289
$ws.Range("G4").Value = 351
$ws.Range("H4").Value = 224
$ws.Range("I4").Value = 192
$ws.Range("J4").Value = 32
$ws.Range("K4").Value = 5502
$ws.Range("L4").Value = 4252
$ws.Range("M4").Value = 1250
$ws.Range("N4").Value = 839
$ws.Range("O4").Value = 411
$ws.Range("P4").Value = 570
$ws.Range("Q4").Value = 184
$ws.Range("R4").Value = -661
$ws.Range("S4").Value = 703
$ws.Range("T4").Value = 909
$ws.Range("U4").Value = -725
$ws.Range("V4").Value = 2903
$ws.Range("W4").Value = 8.630000000000001
$ws.Range("X4").Value = 6.68
$ws.Range("Y4").Value = 26.28
$ws.Range("Z4").Value = 4.45
$ws.Range("AA4").Value = 340.09
$ws.Range("AB4").Value = 30.64
$ws.Range("AC4").Value = 1591
$ws.Range("AD4").Value = 23.1
$ws.Range("AE4").Value = 7031
$ws.Range("AF4").Value = 5.23
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AJ4").Value = 12088268
$ws.Range("AI4").ClearContents()

# Row 5
$ws.Range("D5").Value = 4203
$ws.Range("E5").Value = 333
$ws.Range("F5").Value = 333
$ws.Range("G5").Value = 244
$ws.Range("H5").Value = 109
$ws.Range("I5").Value = 64
$ws.Range("J5").Value = 44
$ws.Range("K5").Value = 6619
$ws.Range("L5").Value = 4262
$ws.Range("M5").Value = 2357
$ws.Range("N5").Value = 1584
$ws.Range("O5").Value = 773
$ws.Range("P5").Value = 570
$ws.Range("Q5").Value = -35
$ws.Range("R5").Value = -752
$ws.Range("S5").Value = 728
$ws.Range("T5").Value = 329
$ws.Range("U5").Value = -364
$ws.Range("V5").Value = 2358
$ws.Range("W5").Value = 7.93
$ws.Range("X5").Value = 2.59
$ws.Range("Y5").Value = 5.31
$ws.Range("Z5").Value = 1.8
$ws.Range("AA5").Value = 180.8
$ws.Range("AB5").Value = 164.07
$ws.Range("AC5").Value = 532
$ws.Range("AD5").Value = 97.5
$ws.Range("AE5").Value = 13270
$ws.Range("AF5").Value = 3.91
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AJ5").Value = 12088268
$ws.Range("AI5").ClearContents()

# Row 6
$ws.Range("D6").Value = 5113
$ws.Range("E6").Value = 347
$ws.Range("F6").Value = 347
$ws.Range("G6").Value = 336
$ws.Range("H6").Value = 256
$ws.Range("I6").Value = 184
$ws.Range("K6").Value = 8212
$ws.Range("L6").Value = 4003
$ws.Range("M6").Value = 4209
$ws.Range("N6").Value = 3370
$ws.Range("P6").Value = 720
$ws.Range("Q6").Value = 248
$ws.Range("R6").Value = -478
$ws.Range("S6").Value = 1152
$ws.Range("T6").Value = 288
$ws.Range("U6").Value = -40
$ws.Range("V6").Value = 1997
$ws.Range("W6").Value = 6.8
$ws.Range("X6").Value = 5
$ws.Range("Y6").Value = 7.45
$ws.Range("Z6").Value = 3.45
$ws.Range("AA6").Value = 95.09
$ws.Range("AB6").Value = 352.32
$ws.Range("AC6").Value = 1448
$ws.Range("AD6").Value = 30.8
$ws.Range("AE6").Value = 23637
$ws.Range("AF6").Value = 1.89
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 14406898
$ws.Range("AI6").ClearContents()

# Row 7
$ws.Range("D7").Value = 5468
$ws.Range("E7").Value = 549
$ws.Range("G7").Value = 455
$ws.Range("H7").Value = 344
$ws.Range("I7").Value = 212
$ws.Range("K7").Value = 9960
$ws.Range("L7").Value = 5428
$ws.Range("M7").Value = 4532
$ws.Range("N7").Value = 3579
$ws.Range("P7").Value = 720
$ws.Range("Q7").Value = 590
$ws.Range("R7").Value = -462
$ws.Range("S7").Value = -179
$ws.Range("T7").Value = 522
$ws.Range("U7").Value = 101
$ws.Range("W7").Value = 10.04
$ws.Range("X7").Value = 6.29
$ws.Range("Y7").Value = 6.1
$ws.Range("Z7").Value = 3.79
$ws.Range("AA7").Value = 119.78
$ws.Range("AC7").Value = 1471
$ws.Range("AD7").Value = 26.65
$ws.Range("AE7").Value = 25116
$ws.Range("AF7").Value = 1.56
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 6247
$ws.Range("E8").Value = 665
$ws.Range("G8").Value = 606
$ws.Range("H8").Value = 461
$ws.Range("I8").Value = 313
$ws.Range("K8").Value = 10506
$ws.Range("L8").Value = 5564
$ws.Range("M8").Value = 4941
$ws.Range("N8").Value = 3897
$ws.Range("P8").Value = 720
$ws.Range("Q8").Value = 631
$ws.Range("R8").Value = -424
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 178
$ws.Range("U8").Value = 440
$ws.Range("W8").Value = 10.65
$ws.Range("X8").Value = 7.38
$ws.Range("Y8").Value = 8.369999999999999
$ws.Range("Z8").Value = 4.5
$ws.Range("AA8").Value = 112.61
$ws.Range("AC8").Value = 2172
$ws.Range("AD8").Value = 18.05
$ws.Range("AE8").Value = 27347
$ws.Range("AF8").Value = 1.43
$ws.Range("AG8").Value = 100
$ws.Range("AH8").Value = 0.26
$ws.Range("AI8").Value = 4.61

# Row 9
$ws.Range("D9").Value = 6808
$ws.Range("E9").Value = 754
$ws.Range("G9").Value = 690
$ws.Range("H9").Value = 523
$ws.Range("I9").Value = 358
$ws.Range("K9").Value = 10838
$ws.Range("L9").Value = 5460
$ws.Range("M9").Value = 5379
$ws.Range("N9").Value = 4222
$ws.Range("P9").Value = 720
$ws.Range("Q9").Value = 767
$ws.Range("R9").Value = -452
$ws.Range("S9").Value = -53
$ws.Range("T9").Value = 180
$ws.Range("U9").Value = 525
$ws.Range("W9").Value = 11.08
$ws.Range("X9").Value = 7.69
$ws.Range("Y9").Value = 8.82
$ws.Range("Z9").Value = 4.9
$ws.Range("AA9").Value = 101.51
$ws.Range("AC9").Value = 2485
$ws.Range("AD9").Value = 15.78
$ws.Range("AE9").Value = 29631
$ws.Range("AF9").Value = 1.32
$ws.Range("AG9").Value = 133
$ws.Range("AH9").Value = 0.34
$ws.Range("AI9").Value = 5.37
